$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.816.46"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.639.68"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'216.22"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.257"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "'0.0636"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'19.71"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "1.865.75"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "1.637.75"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'63.08"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "25.857.19"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D21").Value = "'192.76"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "'9.97"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "'6.35"
$ws.Range("E23").Value = "  +3.00%  "
$ws.Range("E24").Value = "  +4.65%  "
$ws.Range("D26").Value = "'141.88"
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'0.0493"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "'3.24"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "'1.57"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "'0.907"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").Value = "1.133.02"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "'0.545"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "'0.0156"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "'100.72"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").Value = "'0.807"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "1.775.14"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").Value = "'55.41"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'1.44"
$ws.Range("E48").Value = "  +5.89%  "
$ws.Range("D49").Value = "'0.417"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "'2.31"
$ws.Range("E51").Value = "  -0.32%  "
